$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Range("B32").Value = "ул. им. Дзержинского, д. 101"
$ws1.Range("B40").Value = "ул. им. Дзержинского, д. 102"
$ws1.Range("B3").Value = "ул. Горького, д. 128"

$ws1.Activate()
$ws1.Range("K20").Select()
$ws1.Columns.Item(2).ColumnWidth = 55.6640625
